$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column T is column 20 ("Caption"). Row 1 is the header and is left as-is.
# Rows 2 through 307 hold the caption text that needs to be lower-cased.
$lastRow = $ws.Cells.Item(1, 1).End(4).Row
if ($lastRow -lt 307) { $lastRow = 307 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 20)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = $val.ToLower()
    }
}
